$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2428.4285
$ws.Range("I106").Value = 2428.4285
$ws.Range("K106").Value = 2428.4285
$ws.Range("M106").Value = -1797.4285
$ws.Range("H133").Value = 61540.23
$ws.Range("J133").Value = 61540.23
$ws.Range("L133").Value = 61540.23
$ws.Range("N133").Value = -71660.23000000001
$ws.Range("H135").Value = 770513.2
$ws.Range("I135").Value = 1053697.4
$ws.Range("J135").Value = 1870.2858
$ws.Range("K135").Value = 9483276.6
$ws.Range("L135").Value = 16832.5722
$ws.Range("M135").Value = -9480741.6
$ws.Range("N135").Value = -21902.5722
$ws.Range("H137").Value = 1545416
$ws.Range("I137").Value = 1433786
$ws.Range("J137").Value = 1675651
$ws.Range("K137").Value = 4301358
$ws.Range("L137").Value = 5026953
$ws.Range("M137").Value = -4298808
$ws.Range("N137").Value = -5032053
$ws.Range("H138").Value = 5655.8535
$ws.Range("I138").Value = 3242.0435
$ws.Range("J138").Value = 6723.5
$ws.Range("K138").Value = 9726.130500000001
$ws.Range("L138").Value = 20170.5
$ws.Range("M138").Value = -4586.130500000001
$ws.Range("N138").Value = -30450.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2654.277
$ws.Range("I32").Value = 1744.5062
$ws.Range("K32").Value = 1744.5062
$ws.Range("M32").Value = -1457.5062
$ws.Range("H61").Value = 4847
$ws.Range("I61").Value = 4571.5
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 4571.5
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -4359.5
$ws.Range("N61").Value = -6924
$ws.Range("H74").Value = 1335.0193
$ws.Range("I74").Value = 1098.8723
$ws.Range("J74").Value = 3554.8
$ws.Range("K74").Value = 1098.8723
$ws.Range("L74").Value = 3554.8
$ws.Range("M74").Value = -224.8723
$ws.Range("N74").Value = -5302.8
$ws.Range("H77").Value = 1335.0193
$ws.Range("I77").Value = 1098.8723
$ws.Range("J77").Value = 3554.8
$ws.Range("K77").Value = 5494.3615
$ws.Range("L77").Value = 17774
$ws.Range("M77").Value = -1126.3615
$ws.Range("N77").Value = -26510
$ws.Range("H122").Value = 4772.364
$ws.Range("I122").Value = 3362.9092
$ws.Range("J122").Value = 6181.8184
$ws.Range("K122").Value = 10088.7276
$ws.Range("L122").Value = 18545.4552
$ws.Range("M122").Value = -7638.7276
$ws.Range("N122").Value = -23445.4552
$ws.Range("H136").Value = 4847
$ws.Range("I136").Value = 4571.5
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 13714.5
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -11164.5
$ws.Range("N136").Value = -24600

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 27000000
$ws.Range("I4").Value = 500000.5
$ws.Range("J4").Value = 80000000
$ws.Range("K4").Value = 500000.5
$ws.Range("L4").Value = 80000000
$ws.Range("M4").Value = -499888.5
$ws.Range("N4").Value = -80000224
$ws.Range("H7").Value = 411.17648
$ws.Range("I7").Value = 388.1154
$ws.Range("J7").Value = 486.125
$ws.Range("K7").Value = 388.1154
$ws.Range("L7").Value = 486.125
$ws.Range("M7").Value = -275.1154
$ws.Range("N7").Value = -712.125
$ws.Range("H16").Value = 5413.4443
$ws.Range("I16").Value = 1005
$ws.Range("J16").Value = 16875.4
$ws.Range("K16").Value = 1005
$ws.Range("L16").Value = 16875.4
$ws.Range("M16").Value = -718
$ws.Range("N16").Value = -17449.4
$ws.Range("H31").Value = 88113.336
$ws.Range("I31").Value = 2557.6667
$ws.Range("J31").Value = 173669
$ws.Range("K31").Value = 2557.6667
$ws.Range("L31").Value = 173669
$ws.Range("M31").Value = -2262.6667
$ws.Range("N31").Value = -174259
$ws.Range("H34").Value = 88113.336
$ws.Range("I34").Value = 2557.6667
$ws.Range("J34").Value = 173669
$ws.Range("K34").Value = 2557.6667
$ws.Range("L34").Value = 173669
$ws.Range("M34").Value = -2355.6667
$ws.Range("N34").Value = -174073
$ws.Range("H105").Value = 4379.6
$ws.Range("H113").Value = 5413.4443
$ws.Range("I113").Value = 1005
$ws.Range("J113").Value = 16875.4
$ws.Range("K113").Value = 1005
$ws.Range("L113").Value = 16875.4
$ws.Range("M113").Value = 1165
$ws.Range("N113").Value = -21215.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1244.2
$ws.Range("I23").Value = 198
$ws.Range("J23").Value = 1360.4445
$ws.Range("K23").Value = 594
$ws.Range("L23").Value = 4081.3335
$ws.Range("M23").Value = -359
$ws.Range("N23").Value = -4551.333500000001
$ws.Range("H107").Value = 37253.9
$ws.Range("I107").Value = 1669.909
$ws.Range("J107").Value = 58999.668
$ws.Range("K107").Value = 5009.727000000001
$ws.Range("L107").Value = 176999.004
$ws.Range("M107").Value = -3089.727000000001
$ws.Range("N107").Value = -180839.004
$ws.Range("H113").Value = 2180205
$ws.Range("I113").Value = 12346778
$ws.Range("J113").Value = 1653.7142
$ws.Range("K113").Value = 37040334
$ws.Range("L113").Value = 4961.142599999999
$ws.Range("M113").Value = -37038164
$ws.Range("N113").Value = -9301.142599999999
$ws.Range("H132").Value = 503804.7
$ws.Range("I132").Value = 93518.27
$ws.Range("J132").Value = 914091.0600000001
$ws.Range("K132").Value = 841664.4300000001
$ws.Range("L132").Value = 8226819.540000001
$ws.Range("M132").Value = -839134.4300000001
$ws.Range("N132").Value = -8231879.540000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 357
$ws.Range("I2").Value = 402.35715
$ws.Range("K2").Value = 402.35715
$ws.Range("M2").Value = -289.35715
$ws.Range("H97").Value = 941.875
$ws.Range("I97").Value = 739.93335
$ws.Range("J97").Value = 1278.4445
$ws.Range("K97").Value = 739.93335
$ws.Range("L97").Value = 1278.4445
$ws.Range("M97").Value = -243.93335
$ws.Range("N97").Value = -2270.4445
$ws.Range("H107").Value = 1017.1739
$ws.Range("I107").Value = 1152.9166
$ws.Range("J107").Value = 869.0909
$ws.Range("K107").Value = 1152.9166
$ws.Range("L107").Value = 869.0909
$ws.Range("M107").Value = 767.0834
$ws.Range("N107").Value = -4709.0909
$ws.Range("H113").Value = 308211.66
$ws.Range("I113").Value = 477855.56
$ws.Range("J113").Value = 11334.833
$ws.Range("K113").Value = 477855.56
$ws.Range("L113").Value = 11334.833
$ws.Range("M113").Value = -475685.56
$ws.Range("N113").Value = -15674.833
$ws.Range("H122").Value = 8491.200000000001
$ws.Range("I122").Value = 6676.5
$ws.Range("K122").Value = 20029.5
$ws.Range("M122").Value = -17579.5
$ws.Range("H132").Value = 264344.06
$ws.Range("I132").Value = 280345.06
$ws.Range("J132").Value = 168338
$ws.Range("K132").Value = 841035.1799999999
$ws.Range("L132").Value = 505014
$ws.Range("M132").Value = -838505.1799999999
$ws.Range("N132").Value = -510074

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 103007.5
$ws.Range("I40").Value = 114044.445
$ws.Range("K40").Value = 114044.445
$ws.Range("M40").Value = -113908.445
$ws.Range("H43").Value = 1224999.6
$ws.Range("J43").Value = 1516666.5
$ws.Range("L43").Value = 1516666.5
$ws.Range("N43").Value = -1517052.5
$ws.Range("H82").Value = 688.64703
$ws.Range("I82").Value = 693.86664
$ws.Range("K82").Value = 693.86664
$ws.Range("M82").Value = -332.86664
$ws.Range("H85").Value = 688.64703
$ws.Range("I85").Value = 693.86664
$ws.Range("K85").Value = 693.86664
$ws.Range("M85").Value = 554.13336
$ws.Range("H132").Value = 2684.3438
$ws.Range("I132").Value = 1585.4584
$ws.Range("K132").Value = 4756.3752
$ws.Range("M132").Value = -2226.3752
$ws.Range("H136").Value = 387178.6
$ws.Range("I136").Value = 410230.28
$ws.Range("J136").Value = 10667.667
$ws.Range("K136").Value = 1230690.84
$ws.Range("L136").Value = 32003.001
$ws.Range("M136").Value = -1228140.84
$ws.Range("N136").Value = -37103.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H113").Value = 634.86664
$ws.Range("I113").Value = 570.5
$ws.Range("K113").Value = 1711.5
$ws.Range("M113").Value = 458.5
$ws.Range("H132").Value = 24963.045
$ws.Range("I132").Value = 1341.5
$ws.Range("J132").Value = 174566.17
$ws.Range("K132").Value = 4024.5
$ws.Range("L132").Value = 523698.51
$ws.Range("M132").Value = -1494.5
$ws.Range("N132").Value = -528758.51
$ws.Range("H136").Value = 304286.16
$ws.Range("I136").Value = 326395.8
$ws.Range("J136").Value = 225674
$ws.Range("K136").Value = 979187.3999999999
$ws.Range("L136").Value = 677022
$ws.Range("M136").Value = -976637.3999999999
$ws.Range("N136").Value = -682122
$ws.Range("N43").ClearContents()
